# Insert a new data row at row 73 (pushing existing rows 73-91 down to 74-92)
# and populate the newly inserted row with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a whole new row above the current row 73; Excel shifts the
# existing rows (and their formatting/values) down automatically, which
# reproduces rows 74-92 being identical to the old rows 73-91.
$ws.Rows("73").Insert()

# Populate the newly inserted row 73 with the new record.
$ws.Cells.Item(73, 1).Value  = 5
$ws.Cells.Item(73, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(73, 3).Value  = "Maule"
$ws.Cells.Item(73, 4).Value  = 44588
$ws.Cells.Item(73, 5).Value  = 7
$ws.Cells.Item(73, 6).Value  = 100112001
$ws.Cells.Item(73, 7).Value  = "Berenjena"
$ws.Cells.Item(73, 8).Value  = "Sin especificar"
$ws.Cells.Item(73, 9).Value  = "Primera"
$ws.Cells.Item(73, 10).Value = 200
$ws.Cells.Item(73, 11).Value = 7000
$ws.Cells.Item(73, 12).Value = 7000
$ws.Cells.Item(73, 13).Value = 7000
$ws.Cells.Item(73, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(73, 15).Value = "Región del Maule"
$ws.Cells.Item(73, 16).Value = 140
$ws.Cells.Item(73, 17).Value = 50
$ws.Cells.Item(73, 18).Value = "Hortaliza"
